$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 984
$ws.Range("J17").Value = 1070.4286
$ws.Range("L17").Value = 3211.2858
$ws.Range("N17").Value = -3547.2858

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 576.4
$ws.Range("I41").Value = 334.2857
$ws.Range("J41").Value = 706.7692
$ws.Range("K41").Value = 334.2857
$ws.Range("L41").Value = 706.7692
$ws.Range("M41").Value = 105.7143
$ws.Range("N41").Value = -1586.7692

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1878.3636
$ws.Range("I100").Value = 1926
$ws.Range("K100").Value = 1926
$ws.Range("M100").Value = -1385

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 520.1875
$ws.Range("J121").Value = 410.6154
$ws.Range("L121").Value = 1231.8462
$ws.Range("N121").Value = -4725.8462

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1905.6786
$ws.Range("I129").Value = 729.2857
$ws.Range("K129").Value = 2187.8571
$ws.Range("M129").Value = 2812.1429

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3742.6086
$ws.Range("I131").Value = 747.1429000000001
$ws.Range("J131").Value = 5053.125
$ws.Range("K131").Value = 2241.4287
$ws.Range("L131").Value = 15159.375
$ws.Range("M131").Value = 2798.5713
$ws.Range("N131").Value = -25239.375

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1939.0625
$ws.Range("I137").Value = 3116.8572
$ws.Range("J137").Value = 1454.0883
$ws.Range("K137").Value = 9350.571599999999
$ws.Range("L137").Value = 4362.2649
$ws.Range("M137").Value = -6800.571599999999
$ws.Range("N137").Value = -9462.2649

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 12021.667
$ws.Range("I141").Value = 2394.375
$ws.Range("J141").Value = 31276.25
$ws.Range("K141").Value = 7183.125
$ws.Range("L141").Value = 93828.75
$ws.Range("M141").Value = -2003.125
$ws.Range("N141").Value = -104188.75

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1874.9131
$ws.Range("I61").Value = 1873.2174
$ws.Range("J61").Value = 1876.6086
$ws.Range("K61").Value = 1873.2174
$ws.Range("L61").Value = 1876.6086
$ws.Range("M61").Value = -1661.2174
$ws.Range("N61").Value = -2300.6086

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1645.9231
$ws.Range("J74").Value = 2427.2727
$ws.Range("L74").Value = 2427.2727
$ws.Range("N74").Value = -4175.2727

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1645.9231
$ws.Range("J77").Value = 2427.2727
$ws.Range("L77").Value = 12136.3635
$ws.Range("N77").Value = -20872.3635

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1225.3529
$ws.Range("I97").Value = 902.2222
$ws.Range("J97").Value = 1588.875
$ws.Range("K97").Value = 902.2222
$ws.Range("L97").Value = 1588.875
$ws.Range("M97").Value = -406.2222
$ws.Range("N97").Value = -2580.875

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4499.1616
$ws.Range("I132").Value = 3635.0732
$ws.Range("J132").Value = 5811.2964
$ws.Range("K132").Value = 10905.2196
$ws.Range("L132").Value = 17433.8892
$ws.Range("M132").Value = -8375.2196
$ws.Range("N132").Value = -22493.8892

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1874.9131
$ws.Range("I136").Value = 1873.2174
$ws.Range("J136").Value = 1876.6086
$ws.Range("K136").Value = 5619.6522
$ws.Range("L136").Value = 5629.825800000001
$ws.Range("M136").Value = -3069.6522
$ws.Range("N136").Value = -10729.8258

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1246.2188
$ws.Range("I94").Value = 1396.45
$ws.Range("K94").Value = 1396.45
$ws.Range("M94").Value = -945.45

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2925.5557
$ws.Range("I99").Value = 3296.6667
$ws.Range("K99").Value = 3296.6667
$ws.Range("M99").Value = -1798.6667

# BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 10595.75
$ws.Range("J100").Value = 10595.75
$ws.Range("L100").Value = 10595.75
$ws.Range("N100").Value = -12759.75

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2076.8635
$ws.Range("I134").Value = 1307.5778
$ws.Range("J134").Value = 3725.3333
$ws.Range("K134").Value = 3922.7334
$ws.Range("L134").Value = 11175.9999
$ws.Range("M134").Value = -1387.7334
$ws.Range("N134").Value = -16245.9999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2789.761
$ws.Range("I31").Value = 1398.2727
$ws.Range("J31").Value = 4065.2917
$ws.Range("K31").Value = 1398.2727
$ws.Range("L31").Value = 4065.2917
$ws.Range("M31").Value = -1103.2727
$ws.Range("N31").Value = -4655.2917

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2789.761
$ws.Range("I34").Value = 1398.2727
$ws.Range("J34").Value = 4065.2917
$ws.Range("K34").Value = 1398.2727
$ws.Range("L34").Value = 4065.2917
$ws.Range("M34").Value = -1196.2727
$ws.Range("N34").Value = -4469.2917

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2922.6304
$ws.Range("I58").Value = 1290.9667
$ws.Range("J58").Value = 5982
$ws.Range("K58").Value = 1290.9667
$ws.Range("L58").Value = 5982
$ws.Range("M58").Value = -1087.9667
$ws.Range("N58").Value = -6388

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2922.6304
$ws.Range("I136").Value = 1290.9667
$ws.Range("J136").Value = 5982
$ws.Range("K136").Value = 3872.9001
$ws.Range("L136").Value = 17946
$ws.Range("M136").Value = -1322.9001
$ws.Range("N136").Value = -23046

# CUL row 104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 50000250
$ws.Range("I104").Value = 50000250
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 150000750
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -149998129
$ws.Range("N104").ClearContents()

# CUL row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 6504.15
$ws.Range("I124").Value = 2710
$ws.Range("J124").Value = 7768.8667
$ws.Range("K124").Value = 8130
$ws.Range("L124").Value = 23306.6001
$ws.Range("M124").Value = -3220
$ws.Range("N124").Value = -33126.6001

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 6473.75
$ws.Range("I125").Value = 3930
$ws.Range("J125").Value = 8000
$ws.Range("K125").Value = 11790
$ws.Range("L125").Value = 24000
$ws.Range("M125").Value = -6870
$ws.Range("N125").Value = -33840

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1366.2162
$ws.Range("I131").Value = 2175
$ws.Range("J131").Value = 1209.6774
$ws.Range("K131").Value = 6525
$ws.Range("L131").Value = 3629.0322
$ws.Range("M131").Value = -1485
$ws.Range("N131").Value = -13709.0322

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1138.4615
$ws.Range("I97").Value = 787.5
$ws.Range("J97").Value = 1700
$ws.Range("K97").Value = 787.5
$ws.Range("L97").Value = 1700
$ws.Range("M97").Value = -291.5
$ws.Range("N97").Value = -2692

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 542.9474
$ws.Range("I107").Value = 542.9474
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 542.9474
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1377.0526
$ws.Range("N107").ClearContents()

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 477.5
$ws.Range("I22").Value = 448.57144
$ws.Range("J22").Value = 518
$ws.Range("K22").Value = 448.57144
$ws.Range("L22").Value = 518
$ws.Range("M22").Value = -153.57144
$ws.Range("N22").Value = -1108

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 477.5
$ws.Range("I27").Value = 448.57144
$ws.Range("J27").Value = 518
$ws.Range("K27").Value = 448.57144
$ws.Range("L27").Value = 518
$ws.Range("M27").Value = -341.57144
$ws.Range("N27").Value = -732

# LTW row 30
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 33707.668
$ws.Range("I30").Value = 552.5
$ws.Range("J30").Value = 100018
$ws.Range("K30").Value = 552.5
$ws.Range("L30").Value = 100018
$ws.Range("M30").Value = -444.5
$ws.Range("N30").Value = -100234

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41670440
$ws.Range("I40").Value = 76926610
$ws.Range("J40").Value = 4065.3635
$ws.Range("K40").Value = 76926610
$ws.Range("L40").Value = 4065.3635
$ws.Range("M40").Value = -76926474
$ws.Range("N40").Value = -4337.363499999999

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 769810
$ws.Range("J46").Value = 1250616.2
$ws.Range("L46").Value = 1250616.2
$ws.Range("N46").Value = -1250992.2

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1780.6
$ws.Range("I93").Value = 1467.6666
$ws.Range("J93").Value = 2250
$ws.Range("K93").Value = 1467.6666
$ws.Range("L93").Value = 2250
$ws.Range("M93").Value = -219.6666
$ws.Range("N93").Value = -4746

# LTW row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 22988.166
$ws.Range("J94").Value = 22988.166
$ws.Range("L94").Value = 22988.166
$ws.Range("N94").Value = -24340.166

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2422.0476
$ws.Range("I122").Value = 1910.8667
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 5732.6001
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -3282.6001
$ws.Range("N122").Value = -16000
